$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.992.38"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "2.446.61"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'570.45"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "'146.72"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.448.91"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").Value = "'5.23"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "'26.86"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "'0.0000180"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "2.898.18"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "62.884.75"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "2.451.40"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("E20").Value = "  +6.20%  "
$ws.Range("D21").Value = "'324.37"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Value = "'4.17"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'1.94"
$ws.Range("E23").Value = "  +11.95%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'66.26"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "'618.92"
$ws.Range("E26").Value = "  +11.01%  "
$ws.Range("D27").Value = "'8.63"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "'0.0000103"
$ws.Range("E28").Value = "  +10.17%  "
$ws.Range("D29").Value = "2.573.45"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'1.48"
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").Value = "'8.26"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "'0.144"
$ws.Range("E33").Value = "  -2.37%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "'5.10"
$ws.Range("E35").Value = "  +7.93%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.382"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "'18.66"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "'145.70"
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  +15.69%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "'147.17"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'3.72"
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("D47").Value = "'0.0540"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").Value = "'20.61"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "'0.0921"
$ws.Range("E51").Value = "  +0.01%  "
